$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1548
    $ws.Range("F9").Value = 390
}
